$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the header cells in row 1 (columns D, E, F change text; the
# underlying shared-string table is renumbered as a result, but the
# column layout itself is unchanged).
$ws.Range("D1").Value = "下次調整日期"
$ws.Range("E1").Value = "首次調整日期"
$ws.Range("F1").Value = "商品代碼"

# Move the active selection to H11 (matches the saved sheetView selection).
$ws.Range("H11").Select()
